$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Update the wait time value (seconds -> milliseconds)
$ws.Range("B2").Value = 20000

# Update the description text to reflect the new unit
$ws.Range("C2").Value = "Amount of time in millisesondsbot has to wait for user input before proceeding."

# Update the active cell selection
$ws.Range("C6").Select()
